$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Update the "last updated" date in C1 (About sheet) from 3/11/2022 (44631)
# to 9/2/2022 (44806), keeping the existing date number format (style s="4").
$ws.Range("C1").Value = 44806
